$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1797.4517
$ws.Range("I17").Value = 697
$ws.Range("K17").Value = 2091
$ws.Range("M17").Value = -1923

$ws.Range("H33").Value = 144.28572
$ws.Range("I33").Value = 78.46154
$ws.Range("J33").Value = 1000
$ws.Range("K33").Value = 78.46154
$ws.Range("L33").Value = 1000
$ws.Range("M33").Value = 150.53846
$ws.Range("N33").Value = -1458

$ws.Range("H40").Value = 5386.722
$ws.Range("I40").Value = 4050
$ws.Range("J40").Value = 7487.2856
$ws.Range("K40").Value = 4050
$ws.Range("L40").Value = 7487.2856
$ws.Range("M40").Value = -3875
$ws.Range("N40").Value = -7837.2856

$ws.Range("H62").Value = 9250
$ws.Range("I62").Value = 1000
$ws.Range("K62").Value = 1000
$ws.Range("M62").Value = -376

$ws.Range("H65").Value = 9250
$ws.Range("I65").Value = 1000
$ws.Range("K65").Value = 5000
$ws.Range("M65").Value = -1880

$ws.Range("H113").Value = 1987.6428
$ws.Range("I113").Value = 1843.9166
$ws.Range("K113").Value = 1843.9166
$ws.Range("M113").Value = 1410.0834

$ws.Range("H117").Value = 51871
$ws.Range("J117").Value = 51871
$ws.Range("L117").Value = 51871
$ws.Range("N117").Value = -61049

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 847.3
$ws.Range("I2").Value = 847.3
$ws.Range("K2").Value = 847.3
$ws.Range("M2").Value = -734.3

$ws.Range("H10").Value = 1000
$ws.Range("J10").Value = 1000
$ws.Range("L10").Value = 1000
$ws.Range("N10").Value = -1340

$ws.Range("H45").Value = 4023
$ws.Range("I45").Value = 2222
$ws.Range("J45").Value = 4923.5
$ws.Range("K45").Value = 2222
$ws.Range("L45").Value = 4923.5
$ws.Range("M45").Value = -1845
$ws.Range("N45").Value = -5677.5

$ws.Range("H55").Value = 47999.668
$ws.Range("J55").Value = 47999.668
$ws.Range("L55").Value = 47999.668
$ws.Range("N55").Value = -48629.668

$ws.Range("H76").Value = 16666.666
$ws.Range("J76").Value = 16666.666
$ws.Range("L76").Value = 16666.666
$ws.Range("N76").Value = -17342.666

$ws.Range("H79").Value = 16666.666
$ws.Range("J79").Value = 16666.666
$ws.Range("L79").Value = 16666.666
$ws.Range("N79").Value = -19006.666

$ws.Range("H116").Value = 847.3
$ws.Range("I116").Value = 847.3
$ws.Range("K116").Value = 847.3
$ws.Range("M116").Value = 1446.7

$ws.Range("H141").Value = 7499
$ws.Range("J141").Value = 7499
$ws.Range("L141").Value = 7499
$ws.Range("N141").Value = -17859

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 847.3
$ws.Range("I3").Value = 847.3
$ws.Range("K3").Value = 847.3
$ws.Range("M3").Value = -733.3

$ws.Range("H76").Value = 50000
$ws.Range("J76").Value = 50000
$ws.Range("L76").Value = 50000
$ws.Range("N76").Value = -50630

$ws.Range("H79").Value = 50000
$ws.Range("J79").Value = 50000
$ws.Range("L79").Value = 50000
$ws.Range("N79").Value = -52184

$ws.Range("H134").Value = 8983.625
$ws.Range("I134").Value = 1752.7142
$ws.Range("K134").Value = 5258.142599999999
$ws.Range("M134").Value = -2723.142599999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H5").Value = 837.06665
$ws.Range("I5").Value = 180.7
$ws.Range("J5").Value = 2149.8
$ws.Range("K5").Value = 180.7
$ws.Range("L5").Value = 2149.8
$ws.Range("M5").Value = -68.69999999999999
$ws.Range("N5").Value = -2373.8

$ws.Range("H31").Value = 6693.1934
$ws.Range("I31").Value = 3955.5715
$ws.Range("K31").Value = 3955.5715
$ws.Range("M31").Value = -3660.5715

$ws.Range("H34").Value = 6693.1934
$ws.Range("I34").Value = 3955.5715
$ws.Range("K34").Value = 3955.5715
$ws.Range("M34").Value = -3753.5715

$ws.Range("H39").Value = 0
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 0
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("M39").ClearContents()
$ws.Range("N39").ClearContents()

$ws.Range("H49").Value = 0
$ws.Range("I49").Value = 0
$ws.Range("J49").Value = 0
$ws.Range("K49").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("M49").ClearContents()
$ws.Range("N49").ClearContents()

$ws.Range("H140").Value = 137890
$ws.Range("J140").Value = 137890
$ws.Range("L140").Value = 137890
$ws.Range("N140").Value = -148250

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 1260
$ws.Range("I107").Value = 500
$ws.Range("J107").Value = 1450
$ws.Range("K107").Value = 1500
$ws.Range("L107").Value = 4350
$ws.Range("M107").Value = 420
$ws.Range("N107").Value = -8190

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H12").Value = 2148.7058
$ws.Range("J12").Value = 1087
$ws.Range("L12").Value = 1087
$ws.Range("N12").Value = -1367

$ws.Range("H132").Value = 114338.336
$ws.Range("I132").Value = 127830.625
$ws.Range("K132").Value = 383491.875
$ws.Range("M132").Value = -380961.875

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1670.6364
$ws.Range("I22").Value = 1123.7142
$ws.Range("J22").Value = 2073.6316
$ws.Range("K22").Value = 1123.7142
$ws.Range("L22").Value = 2073.6316
$ws.Range("M22").Value = -828.7141999999999
$ws.Range("N22").Value = -2663.6316

$ws.Range("H26").Value = 3999
$ws.Range("I26").Value = 3999
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 3999
$ws.Range("L26").Value = 0
$ws.Range("M26").Value = -3704
$ws.Range("N26").ClearContents()

$ws.Range("H27").Value = 1670.6364
$ws.Range("I27").Value = 1123.7142
$ws.Range("J27").Value = 2073.6316
$ws.Range("K27").Value = 1123.7142
$ws.Range("L27").Value = 2073.6316
$ws.Range("M27").Value = -1016.7142
$ws.Range("N27").Value = -2287.6316

$ws.Range("H46").Value = 6726.727
$ws.Range("I46").Value = 999
$ws.Range("J46").Value = 7299.5
$ws.Range("K46").Value = 999
$ws.Range("L46").Value = 7299.5
$ws.Range("M46").Value = -811
$ws.Range("N46").Value = -7675.5

$ws.Range("H56").Value = 13137.857
$ws.Range("I56").Value = 11160.833
$ws.Range("J56").Value = 25000
$ws.Range("K56").Value = 11160.833
$ws.Range("L56").Value = 25000
$ws.Range("M56").Value = -10469.833
$ws.Range("N56").Value = -26382

$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("N122").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H13").Value = 1900
$ws.Range("I13").Value = 800
$ws.Range("J13").Value = 2450
$ws.Range("K13").Value = 800
$ws.Range("L13").Value = 2450
$ws.Range("M13").Value = -660
$ws.Range("N13").Value = -2730

$ws.Range("H51").Value = 21401.428
$ws.Range("I51").Value = 21735
$ws.Range("K51").Value = 21735
$ws.Range("M51").Value = -21225

$ws.Range("H82").Value = 33987.5
$ws.Range("J82").Value = 33987.5
$ws.Range("L82").Value = 33987.5
$ws.Range("N82").Value = -34753.5

$ws.Range("H85").Value = 33987.5
$ws.Range("J85").Value = 33987.5
$ws.Range("L85").Value = 33987.5
$ws.Range("N85").Value = -36639.5

$ws.Range("H132").Value = 1499.625
$ws.Range("I132").Value = 1499.625
$ws.Range("K132").Value = 4498.875
$ws.Range("M132").Value = -1968.875
